# Grump metadata edit: add Longhurst provinces and Season variable rows
# to the "vars_meta_data" sheet (rows 41-43), and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")
$ws.Activate()

# --- Carry over the formatting already used by the surrounding template
#     rows (e.g. row 39) onto the new rows' A:D and G columns before
#     writing the new values into them. ---
$ws.Range("B39:D39").Copy()
$ws.Range("A41:C43").PasteSpecial(-4122)
$ws.Range("G39").Copy()
$ws.Range("G41:G43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 41: Longhurst_Long ---
$ws.Cells.Item(41, 1).Value = "Longhurst_Long"
$ws.Cells.Item(41, 2).Value = "Longhurst province sample was taken in."
$ws.Cells.Item(41, 3).Value = "NA"
$ws.Cells.Item(41, 4).Value = "NA"
$ws.Cells.Item(41, 5).Value = "Irregular"
$ws.Cells.Item(41, 6).Value = "Irregular"
$ws.Cells.Item(41, 7).Value = "Biology"
$ws.Cells.Item(41, 8).Value = 1

# --- Row 42: Longhurst_Short ---
$ws.Cells.Item(42, 1).Value = "Longhurst_Short"
$ws.Cells.Item(42, 2).Value = "Longhurst province sample was taken in, shortened code."
$ws.Cells.Item(42, 3).Value = "NA"
$ws.Cells.Item(42, 4).Value = "NA"
$ws.Cells.Item(42, 5).Value = "Irregular"
$ws.Cells.Item(42, 6).Value = "Irregular"
$ws.Cells.Item(42, 7).Value = "Biology"
$ws.Cells.Item(42, 8).Value = 1

# --- Row 43: Season ---
$ws.Cells.Item(43, 1).Value = "Season"
$ws.Cells.Item(43, 2).Value = "Season sample was taken in."
$ws.Cells.Item(43, 3).Value = "NA"
$ws.Cells.Item(43, 4).Value = "NA"
$ws.Cells.Item(43, 5).Value = "Irregular"
$ws.Cells.Item(43, 6).Value = "Irregular"
$ws.Cells.Item(43, 7).Value = "Biology"
$ws.Cells.Item(43, 8).Value = 1

# --- Update the active selection to match the newly entered block ---
$ws.Range("A41:H43").Select()
